# "Update parallel merge - limit."
#
# The second table (originally rows 8-12, describing the "Parellel merge"
# runtime table) is restructured:
#   - row 8 header: B8 changes from "n" to "threads" (shared string index 5
#     is retargeted from "runtime" -> "threads"), and C8:F8 now hold the
#     thread-count series 1, 2, 4, 8 (replacing the old single "runtime"
#     label cell).
#   - a new row is inserted at row 9 holding the "n" label that used to be
#     in B8, pushing the old n-value rows (10/100/1000/10000) down from
#     rows 9-12 to rows 10-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 9; this shifts the old B9:B12 (10,100,1000,10000)
# down to B10:B13 and grows the used range to A1:H13.
$ws.Rows("9:9").Insert()

# Re-label the header row: B8 "n" -> "threads", and populate the new
# thread-count columns C8:F8.
$ws.Range("B8").Value = "threads"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 8

# The freshly inserted row 9 becomes the new "n" label (what used to live
# in B8).
$ws.Range("B9").Value = "n"

# Match the saved selection state (active cell moved from G12 to H12).
$ws.Range("H12").Select()
